# Fruta / hortaliza, semanal
# Insert two new weekly observation rows into the "Pepino ensalada" sheet,
# right after the existing row 590, shifting all subsequent rows down by 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 591 (they inherit formatting
# from the row being pushed down, matching how Excel's Insert normally behaves).
$ws.Rows("591:592").Insert()

# --- New row 591 ---
$ws.Range("A591").Value = 5
$ws.Range("B591").Value = "Macroferia Regional de Talca"
$ws.Range("C591").Value = "Maule"
$ws.Range("D591").Value = 45021
$ws.Range("E591").Value = 7
$ws.Range("F591").Value = 100112043
$ws.Range("G591").Value = "Pepino ensalada"
$ws.Range("H591").Value = "Sin especificar"
$ws.Range("I591").Value = "Primera"
$ws.Range("J591").Value = 300
$ws.Range("K591").Value = 7000
$ws.Range("L591").Value = 7000
$ws.Range("M591").Value = 7000
$ws.Range("N591").Value = "$/caja 60 unidades"
$ws.Range("O591").Value = "Región de Arica y Parinacota"
$ws.Range("P591").Value = 117
$ws.Range("Q591").Value = 60
$ws.Range("R591").Value = "Hortaliza"

# --- New row 592 ---
$ws.Range("A592").Value = 5
$ws.Range("B592").Value = "Macroferia Regional de Talca"
$ws.Range("C592").Value = "Maule"
$ws.Range("D592").Value = 45021
$ws.Range("E592").Value = 7
$ws.Range("F592").Value = 100112043
$ws.Range("G592").Value = "Pepino ensalada"
$ws.Range("H592").Value = "Sin especificar"
$ws.Range("I592").Value = "Primera"
$ws.Range("J592").Value = 200
$ws.Range("K592").Value = 9000
$ws.Range("L592").Value = 9000
$ws.Range("M592").Value = 9000
$ws.Range("N592").Value = "$/caja 80 unidades"
$ws.Range("O592").Value = "Región del Maule"
$ws.Range("P592").Value = 112
$ws.Range("Q592").Value = 80
$ws.Range("R592").Value = "Hortaliza"
